$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RunControl")
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()
$ws.Range("A6:AX6").Copy()
$ws.Range("A13:AX13").PasteSpecial(-4104)
$ws.Range("A12").Value = "check"
